# Helper: assign a value to a cell while forcing Excel to store it as
# literal text (so numeric-looking strings like "310.47" or "1.00" are not
# silently converted into floating point numbers), then restore the cell
# style to the default "Normal" style so no stray number-format is left on
# the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15 and 16: Chainlink and WrappedliquidstakedEther2.0 swap places
Set-TextValue $ws.Range("B15") "Chainlink"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D15") "15.29"
$ws.Range("E15").Value = "  -1.41%  "

Set-TextValue $ws.Range("B16") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D16") "2.581.54"
$ws.Range("E16").Value = "  -3.26%  "

# Remaining rows: update Price (D) and Volume(1h) (E) columns
Set-TextValue $ws.Range("D2") "43.077.16"
$ws.Range("E2").Value = "  +0.75%  "
Set-TextValue $ws.Range("D3") "2.305.64"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws.Range("D5") "310.47"
$ws.Range("E5").Value = "  -2.33%  "
Set-TextValue $ws.Range("D6") "104.78"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.06%  "
Set-TextValue $ws.Range("D9") "0.607"
$ws.Range("E9").Value = "  -0.57%  "
Set-TextValue $ws.Range("D10") "39.61"
$ws.Range("E10").Value = "  -1.46%  "
Set-TextValue $ws.Range("D11") "0.0907"
$ws.Range("E11").Value = "  -0.35%  "
Set-TextValue $ws.Range("D12") "8.26"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("E13").Value = "  +0.02%  "
Set-TextValue $ws.Range("D14") "0.998"
$ws.Range("E14").Value = "  +0.90%  "
Set-TextValue $ws.Range("D17") "2.306.85"
$ws.Range("E17").Value = "  -0.44%  "
Set-TextValue $ws.Range("D18") "42.875.99"
$ws.Range("E18").Value = "  +0.33%  "
Set-TextValue $ws.Range("D19") "7.33"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("E20").Value = "  -1.27%  "
Set-TextValue $ws.Range("D21") "13.59"
$ws.Range("E21").Value = "  +0.65%  "
Set-TextValue $ws.Range("D22") "73.40"
$ws.Range("E22").Value = "  -0.87%  "
Set-TextValue $ws.Range("D23") "3.44"
$ws.Range("E23").Value = "  -2.93%  "
Set-TextValue $ws.Range("D24") "267.34"
$ws.Range("E24").Value = "  -0.98%  "
Set-TextValue $ws.Range("D25") "2.25"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  +0.50%  "
Set-TextValue $ws.Range("D27") "10.98"
$ws.Range("E27").Value = "  +0.26%  "
Set-TextValue $ws.Range("D28") "7.33"
$ws.Range("E28").Value = "  +11.99%  "
$ws.Range("E29").Value = "  -2.58%  "
Set-TextValue $ws.Range("D30") "22.28"
$ws.Range("E30").Value = "  -2.04%  "
Set-TextValue $ws.Range("D31") "37.29"
$ws.Range("E31").Value = "  -3.02%  "
Set-TextValue $ws.Range("D32") "164.85"
$ws.Range("E32").Value = "  -1.26%  "
Set-TextValue $ws.Range("D33") "0.0861"
$ws.Range("E33").Value = "  -3.18%  "
Set-TextValue $ws.Range("D34") "2.87"
$ws.Range("E34").Value = "  +9.23%  "
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  -2.79%  "
Set-TextValue $ws.Range("D37") "4.57"
$ws.Range("E37").Value = "  -0.73%  "
Set-TextValue $ws.Range("D38") "0.0349"
$ws.Range("E38").Value = "  -1.72%  "
Set-TextValue $ws.Range("D39") "2.86"
$ws.Range("E39").Value = "  +2.40%  "
Set-TextValue $ws.Range("D40") "3.62"
$ws.Range("E40").Value = "  -2.97%  "
Set-TextValue $ws.Range("D41") "108.17"
$ws.Range("E41").Value = "  +7.46%  "
Set-TextValue $ws.Range("D42") "1.59"
$ws.Range("E42").Value = "  -3.62%  "
Set-TextValue $ws.Range("D43") "71.53"
$ws.Range("E43").Value = "  +1.13%  "
Set-TextValue $ws.Range("D44") "0.228"
$ws.Range("E44").Value = "  +0.56%  "
Set-TextValue $ws.Range("D45") "1.00"
$ws.Range("E45").Value = "  -0.39%  "
Set-TextValue $ws.Range("D46") "12.25"
$ws.Range("E46").Value = "  -1.08%  "
Set-TextValue $ws.Range("D47") "1.722.04"
$ws.Range("E47").Value = "  +4.49%  "
Set-TextValue $ws.Range("D48") "111.42"
$ws.Range("E48").Value = "  -5.31%  "
Set-TextValue $ws.Range("D49") "76.64"
$ws.Range("E49").Value = "  -6.48%  "
Set-TextValue $ws.Range("D50") "8.73"
$ws.Range("E50").Value = "  -1.77%  "
Set-TextValue $ws.Range("D51") "5.16"
$ws.Range("E51").Value = "  -3.12%  "
